# Members sheet update:
#  - insert two new members ("Lisanne Canjels", "Chris Baeken") into the
#    member table at rows 25-26 (pushing the existing rows 25-33 down to 27-35)
#  - renumber the "nummer" column (B) sequentially for the shifted block
#  - highlight the renumbered "nummer" cells (B23:B35) with a yellow fill
#  - mark a handful of additional "zelf tekstje geschreven?" / "feedback/nieuwe
#    tekst ontvangen?" cells (columns D/E) with "v"
#  - widen column A and move the frozen-pane/selection to the bottom of the list

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert 2 new rows for the new members, shifting rows 25-33 -> 27-35 ---
$ws.Rows("25:26").Insert()

# --- new member rows ---
$ws.Range("A25").Value = "Lisanne Canjels"
$ws.Range("B25").Value = 24
$ws.Range("D25").Value = "v"
$ws.Range("E25").Value = "v"

$ws.Range("A26").Value = "Chris Baeken"
$ws.Range("B26").Value = 25
$ws.Range("C26").Value = "affiliate"
$ws.Range("D26").Value = "v"

# --- renumber "nummer" column for the rows that shifted down (old B 24..32 -> new 26..34) ---
$ws.Range("B27").Value = 26
$ws.Range("B28").Value = 27
$ws.Range("B29").Value = 28
$ws.Range("B30").Value = 29
$ws.Range("B31").Value = 30
$ws.Range("B32").Value = 31
$ws.Range("B33").Value = 32
$ws.Range("B34").Value = 33
$ws.Range("B35").Value = 34

# --- add the newly-checked "v" marks on existing members ---
$ws.Range("E9").Value = "v"
$ws.Range("E11").Value = "v"
$ws.Range("E13").Value = "v"
$ws.Range("E14").Value = "v"
$ws.Range("E17").Value = "v"
$ws.Range("E21").Value = "v"
$ws.Range("E23").Value = "v"
$ws.Range("E28").Value = "v"
$ws.Range("E30").Value = "v"
$ws.Range("E32").Value = "v"

# --- highlight the renumbered "nummer" cells yellow ---
$ws.Range("B23:B35").Interior.Color = 65535

# --- widen column A ---
$ws.Columns(1).ColumnWidth = 13.166666666666666

# --- move the active selection down to the bottom of the list (the sheet stays
#     frozen on row 1; the view simply scrolls so row 20 is the first visible
#     row beneath the frozen header, matching where the new rows now sit) ---
$ws.Range("B35").Select()
